# Update the "想去人数" (F column) numeric values on sheets 展览, 演出, and 全部类型
# to match the regenerated data as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1326
$ws1.Range("F4").Value  = 81
$ws1.Range("F6").Value  = 392
$ws1.Range("F7").Value  = 176
$ws1.Range("F8").Value  = 124
$ws1.Range("F9").Value  = 1014
$ws1.Range("F10").Value = 332
$ws1.Range("F12").Value = 42
$ws1.Range("F16").Value = 756
$ws1.Range("F17").Value = 136
$ws1.Range("F18").Value = 703
$ws1.Range("F19").Value = 254
$ws1.Range("F20").Value = 66
$ws1.Range("F21").Value = 973
$ws1.Range("F23").Value = 245
$ws1.Range("F24").Value = 78
$ws1.Range("F25").Value = 356
$ws1.Range("F27").Value = 33

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 354
$ws2.Range("F5").Value = 34

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1326
$ws4.Range("F6").Value  = 81
$ws4.Range("F8").Value  = 392
$ws4.Range("F9").Value  = 176
$ws4.Range("F10").Value = 124
$ws4.Range("F11").Value = 1014
$ws4.Range("F12").Value = 332
$ws4.Range("F15").Value = 42
$ws4.Range("F16").Value = 354
$ws4.Range("F18").Value = 34
$ws4.Range("F23").Value = 756
$ws4.Range("F24").Value = 136
$ws4.Range("F25").Value = 703
$ws4.Range("F26").Value = 254
$ws4.Range("F27").Value = 66
$ws4.Range("F28").Value = 973
$ws4.Range("F32").Value = 245
$ws4.Range("F33").Value = 78
$ws4.Range("F34").Value = 356
$ws4.Range("F38").Value = 33
